# Sathya Coatings ERP Insights — add a new "P006" project row to the
# "Costing & Budgeting" sheet (mirrors the previous "P005" row, but left
# at "Pending Start"), then leave the UI selection where Excel would
# land after the user finished typing the row (cell I15).

$wb = $excel.ActiveWorkbook

# Look the sheet up by name rather than position so this keeps working
# even if sheets get reordered.
$ws = $wb.Worksheets.Item("Costing & Budgeting")

$ws.Range("A7").Value = "P006"
$ws.Range("B7").Value = 6363
$ws.Range("C7").Value = 34
$ws.Range("D7").Value = 44545
$ws.Range("E7").Value = 65565
$ws.Range("F7").Value = 6556
$ws.Range("G7").Value = "Pending Start"

# Match the saved cursor position recorded in the workbook after the edit.
$ws.Range("I15").Select()
